$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.148.03"
$ws.Range("E2").Value = "  -0.79%  "

$ws.Range("D3").Value = "1.905.38"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").Value = "'0.7233"
$ws.Range("E5").Value = "  -6.05%  "

$ws.Range("D6").Value = "'242.92"
$ws.Range("E6").Value = "  -1.96%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("D8").Value = "'0.3113"
$ws.Range("E8").Value = "  -3.00%  "

$ws.Range("D9").Value = "'26.54"
$ws.Range("E9").Value = "  -5.32%  "

$ws.Range("D10").Value = "'0.06872"
$ws.Range("E10").Value = "  -3.39%  "

$ws.Range("D11").Value = "'0.7731"
$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("D12").Value = "'0.07964"
$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("D13").Value = "1.885.66"
$ws.Range("E13").Value = "  -2.69%  "

$ws.Range("D14").Value = "'5.230"
$ws.Range("E14").Value = "  -2.84%  "

$ws.Range("D15").Value = "'91.07"
$ws.Range("E15").Value = "  -4.29%  "

$ws.Range("D16").Value = "30.125.19"
$ws.Range("E16").Value = "  -0.89%  "

$ws.Range("D17").Value = "'14.12"
$ws.Range("E17").Value = "  -3.13%  "

$ws.Range("D18").Value = "'5.828"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007755"
$ws.Range("E19").Value = "  -3.31%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'238.10"
$ws.Range("E20").Value = "  -7.26%  "

$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "2.143.97"
$ws.Range("E22").Value = "  -2.10%  "

$ws.Range("D23").Value = "'1.005"
$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("D24").Value = "'6.817"
$ws.Range("E24").Value = "  +0.75%  "

$ws.Range("D25").Value = "'9.340"
$ws.Range("E25").Value = "  -2.87%  "

$ws.Range("D26").Value = "'164.85"
$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("D27").Value = "'19.00"
$ws.Range("E27").Value = "  -0.80%  "

$ws.Range("D28").Value = "'0.1262"
$ws.Range("E28").Value = "  -5.94%  "

$ws.Range("D29").Value = "'2.073"
$ws.Range("E29").Value = "  -9.92%  "

$ws.Range("D30").Value = "'1.353"
$ws.Range("E30").Value = "  -1.00%  "

$ws.Range("D31").Value = "'1.542"
$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("D32").Value = "'4.279"
$ws.Range("E32").Value = "  -3.69%  "

$ws.Range("D33").Value = "'4.063"
$ws.Range("E33").Value = "  -2.11%  "

$ws.Range("D34").Value = "'0.05135"
$ws.Range("E34").Value = "  -1.36%  "

$ws.Range("D35").Value = "'1.278"
$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("D36").Value = "'0.7384"
$ws.Range("E36").Value = "  -1.78%  "

$ws.Range("D37").Value = "'2.756"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("D38").Value = "'0.01928"
$ws.Range("E38").Value = "  -2.36%  "

$ws.Range("D39").Value = "'2.793"
$ws.Range("E39").Value = "  -0.57%  "

$ws.Range("D40").Value = "'6.342"
$ws.Range("E40").Value = "  -1.92%  "

$ws.Range("D41").Value = "'74.23"
$ws.Range("E41").Value = "  -6.11%  "

$ws.Range("D42").Value = "'0.4414"
$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("D43").Value = "'1.922"
$ws.Range("E43").Value = "  -2.99%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").Value = "'100.91"
$ws.Range("E46").Value = "  -0.55%  "

$ws.Range("D47").Value = "'7.535"
$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").Value = "'9.686"
$ws.Range("E48").Value = "  -1.20%  "

$ws.Range("D49").Value = "'37.56"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.044.07"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'937.61"
$ws.Range("E51").Value = "  -5.01%  "
